# Weekly update: insert a new reporting week (2022-03-24) at the top of the
# "Comercializadora del Agro de Limarí - Tuna" block (rows 47-63), pushing the
# older weeks down by 3 rows (the block is always 3 quality grades per week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above row 47 - everything from 47 downward shifts to 50+.
$ws.Rows(47).Resize(3).Insert()

# Common (constant-across-rows) column values for this product block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107011
$categoria   = "Tuna"
$variedad    = "Sin especificar"
$unidad      = "$/caja 18 kilos"
$origen      = "Provincia de Limarí"
$kgUnidad    = 18

# NOTE: named parameters don't bind values reliably in this PS host, so
# Set-DataRow is called with strictly positional arguments.
function Set-DataRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New top week: 2022-03-24 (serial 44644), three quality grades.
Set-DataRow 47 44644 "Especial" 300 11500 12000 11750 653
Set-DataRow 48 44644 "Primera"  240 9500  10000 9750  542
Set-DataRow 49 44644 "Segunda"  240 5000  6000  5500  306

# Ensure the Fecha cells use the existing date format (same as every other D-column cell).
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat

Write-Output "done"
